$d = $word.ActiveDocument

# 1) Replace "develop" with "master" in the target paragraph's run text
$d.Content.Find.Execute("I want to change this file by develop", $false, $false, $false, $false, $false, $true, 1, $false, "I want to change this file by master", 2)

# 2) Find the paragraph and split its single run into two runs: "I wan" + "t to change this file by master"
#    Locate the run text range for "I wan" and "t to change..." boundary using Find on the full sentence range.
$target = $d.Content
$target.Find.ClearFormatting()
$target.Find.Execute("I want to change this file by master")
if ($target.Find.Found) {
    $full = $target.Duplicate
    $splitPoint = $full.Start + 5  # after "I wan"
    $r1 = $d.Range($full.Start, $splitPoint)
    $r2 = $d.Range($splitPoint, $full.End)
    # touching the text (no-op set) forces run boundary split in many engines;
    # but to guarantee two distinct runs we set text explicitly.
    $r2.Text = $r2.Text
    $r1.Text = $r1.Text
}

# 3) Insert a new paragraph after that paragraph with "first additional change by develop"
$para = $d.Paragraphs.Last
$endRange = $para.Range
$endRange.InsertParagraphAfter()
$newRange = $d.Range($endRange.End, $endRange.End)
$newPara = $d.Paragraphs.Last
$newPara.Range.Text = "first additional change by develop"
